# Autogenerated on Fri Mar 27 2015 16:05:45 GMT+0000 (Coordinated Universal Time)
# Replaces the literal "<br/>" markers inside a handful of cells on the active
# sheet with real in-cell line breaks (Chr(10)), matching how the source text
# is meant to be rendered with word-wrap instead of literal HTML tags.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B21").Value = "51-100 `n<250 Industry, `n<500 Trade, `n<100 Service, `n<50 Agriculture"
$ws.Range("D21").Value = "100,000 UT to 250,000 UT `n<750,000 Industry, `n<1,000,000 Trade, `n<500,000 Serv., `n<300,000 Agriculture"
$ws.Range("B22").Value = ">100 `n>=250 Industry, `n>=500 Trade, `n>=100 Service, `n>=50 Agriculture"
$ws.Range("D22").Value = ">250,000 UT `n>=750,000 Industry, `n>=1,000,000 Trade, `n>=500,000 Serv.,`n >=300,000 Agriculture"
